$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: prefix with an apostrophe to force text entry (like typing
# into Excel with a leading single-quote), then reset the style to Normal so no
# stray quotePrefix/number-format style sticks to the cell (matches the source
# cells, which carry no explicit style).

$ws.Range("D2").Value = "'63.221.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.47%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.572.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.08%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'585.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.13%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'147.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.69%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +3.31%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +3.80%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E11").Value = "'  +0.12%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.51%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'27.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.17%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.035.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.15%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'63.198.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.45%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +4.00%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.584.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.77%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'11.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.90%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'342.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.07%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +2.83%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.90%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.29%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'66.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.697.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.28%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.24%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +2.21%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +12.04%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'InternetComputer(DFINITY)"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'8.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.67%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'SuiNetwork"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'1.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.79%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'Binance-PegBSC-USD"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.01%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +6.08%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +2.37%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'461.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +13.55%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +3.81%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'176.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.27%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +2.34%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'19.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.37%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'4.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.91%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E40").Value = "'  -0.83%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.05%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'151.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.97%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.05%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'21.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.54%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +6.18%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.615"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.60%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +2.65%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.47%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.31%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.17%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +3.80%  "
$ws.Range("E51").Style = "Normal"
